$d = $word.ActiveDocument

# Locate the "line: 4" paragraph inside the test-log table (it sits right
# after the "header: ..." paragraph in the same table cell) and remove the
# whole paragraph, since this test case's log no longer records a parsed
# line number.
$paras = $d.Paragraphs
for ($i = $paras.Count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "line: 4" -or $t.StartsWith("line: ")) {
        $p.Range.Delete()
        break
    }
}
